$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Real time (minutes)" column D for the newly completed Withdraw tasks
$ws.Range("D26").Value = 10
$ws.Range("D27").Value = 5
$ws.Range("D28").Value = 60

# Move the active selection (as recorded in the saved sheet view) to B34
$ws.Range("B34").Select()
